$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D-column prices that only changed value (no row shuffling)
$ws.Range("D2").Value = "245.58"
$ws.Range("D3").Value = "25.51"
$ws.Range("D4").Value = "5.112"
$ws.Range("D5").Value = "0.05565"
$ws.Range("D6").Value = "6.473"
$ws.Range("D7").Value = "3.018"
$ws.Range("D8").Value = "0.8185"
$ws.Range("D9").Value = "0.8456"

# Rows 10-15 shift up: row10 <- old row11 values, etc, with some price updates
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "0.009748"
$ws.Range("E10").Value = "9OneONEBestin24h"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1340"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "0.03196"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.02855"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09394"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001510"
$ws.Range("E15").Value = "14BitForexTokenBF"

# Simple value-only updates
$ws.Range("D16").Value = "0.006191"
$ws.Range("D18").Value = "2.092"
$ws.Range("D20").Value = "0.06952"
$ws.Range("D22").Value = "3.756"
$ws.Range("D23").Value = "0.04730"
$ws.Range("D25").Value = "0.001247"
$ws.Range("D26").Value = "0.004629"
$ws.Range("D27").Value = "0.00009704"
$ws.Range("E27").Value = "26NitroExNTX"
$ws.Range("D40").Value = "0.03657"

# Rows 41-43 shift, with price updates
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.006222"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1052"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.002501"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "0.008294"
$ws.Range("D45").Value = "0.00005300"
$ws.Range("D47").Value = "0.1331"
$ws.Range("D48").Value = "0.002123"
